$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("{m:template toBulletList(s: String)}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$end = $rng.End

# Step 1: split at offset 11 via no-op toggle -> two runs (both keep rsid for now)
$splitOffset = 11
$beforePublic = $d.Range($start, $start + $splitOffset)
$beforePublic.Font.Bold = 1
$beforePublic.Font.Bold = 0

# Step 2: insert "public " strictly INSIDE run2 (away from its edges) so it merges to ONE
#         run spanning the whole paragraph text while PRESERVING rsid (content edit interior to a run keeps rsid)
$insertAt = $start + $splitOffset + 1
$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.InsertAfter("public ")
$insertedLen = 7
$fullEnd = $end + $insertedLen

# Step 3: mark the first 11 characters (the future run1, "{m:template") with a REAL (non reverted-yet)
#         formatting change so the engine treats it as a genuinely distinct run from the rest
$run1Mark = $d.Range($start, $start + $splitOffset)
$run1Mark.Font.Bold = 1

# Step 4: perform a content edit exactly at the run2/run3 boundary (between " public" and " toBulletList...")
#         this is a structural (edge) edit, so it strips rsid from the touched run-group (run2+run3),
#         but leaves run1 (which now has distinct Bold formatting) alone
$boundary = $start + $splitOffset + $insertedLen
$dummy = $d.Range($boundary, $boundary)
$dummy.InsertAfter("Q")
$dummyDel = $d.Range($boundary, $boundary + 1)
$dummyDel.Delete()

# Step 5: revert run1's Bold mark back to its original (False) value -- this leaves the rPr
#         clean (no stray explicit w:b element) while keeping run1 as its own separate run
$run1Unmark = $d.Range($start, $start + $splitOffset)
$run1Unmark.Font.Bold = 0

# Step 6: split the remaining merged " public toBulletList(s: String)}" run into two runs
#         (" public" and " toBulletList(s: String)}") via a no-op formatting toggle -- this keeps
#         whatever rsid-state currently exists (none) on both resulting pieces
$secondBoundary = $start + $splitOffset + $insertedLen
$thirdPiece = $d.Range($secondBoundary, $fullEnd)
$thirdPiece.Font.Bold = 1
$thirdPiece.Font.Bold = 0
